$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text while we assign the new values,
# so numeric-looking strings (e.g. "244.15") are not auto-converted to numbers.
# The original cells are plain (unstyled) inline strings, so we revert the
# temporary formatting afterwards with ClearFormats().
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Update Price (D) and Volume(1h) (E) columns for unchanged-identity rows ---
$ws.Range("D2").Value = "30.510.11"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.890.83"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "244.15"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4715"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "0.2898"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "0.06498"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "22.27"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").Value = "0.07760"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "1.888.02"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "95.88"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "0.7270"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").Value = "5.195"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "282.66"
$ws.Range("E16").Value = "  +3.45%  "
$ws.Range("D17").Value = "30.490.23"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "13.09"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D20").Value = "0.000007481"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "2.136.45"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D24").Value = "6.320"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("D25").Value = "164.04"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").Value = "9.104"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "18.90"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").Value = "1.895"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "1.335"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").Value = "0.09701"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("D31").Value = "1.470"
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("D32").Value = "4.285"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").Value = "4.151"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").Value = "0.04862"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "1.129"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D36").Value = "0.6949"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").Value = "2.715"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "0.01890"
$ws.Range("E38").Value = "  +1.90%  "
$ws.Range("D39").Value = "2.815"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("D40").Value = "75.24"
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("D41").Value = "6.211"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").Value = "1.985"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").Value = "0.4278"
$ws.Range("E43").Value = "  +2.75%  "
$ws.Range("D45").Value = "0.8296"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "101.41"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "9.647"
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("D48").Value = "6.969"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "35.23"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "909.57"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "0.05753"
$ws.Range("E51").Value = "  +1.69%  "

# --- Row 44: only Volume(1h) changes, Price (D44) stays "1.000" ---
$ws.Range("E44").Value = "  -0.10%  "

# --- Rows 22 and 23 swap coin identity (Uniswap <-> BinanceUSD) with new data ---
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.284"
$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "1.0000"
$ws.Range("E23").Value = "  -0.12%  "

# Revert the temporary text formatting on column D back to the original (unstyled) state.
$priceRange.ClearFormats()
